$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 14.47727272727272
$ws.Range("N2").Value = 1.815485677363773
$ws.Range("O2").Value = 1.979371877230549

$ws.Range("I3").Value = 1.791666666666668
$ws.Range("N3").Value = 1.620655622136059
$ws.Range("O3").Value = 1.747323835194455

$ws.Range("I5").Value = 14.47727272727272
$ws.Range("N5").Value = 1.815485677363773
$ws.Range("O5").Value = 1.979371877230549

$ws.Range("I7").Value = 14.47727272727272
$ws.Range("N7").Value = 1.815485677363773
$ws.Range("O7").Value = 1.979371877230549

$ws.Range("I8").Value = 14.47727272727272
$ws.Range("N8").Value = 1.815485677363773
$ws.Range("O8").Value = 1.979371877230549

$ws.Range("I9").Value = 13.76976495726495
$ws.Range("N9").Value = 1.803394296576035
$ws.Range("O9").Value = 1.964819060413116

$ws.Range("I11").Value = 1.791666666666668
$ws.Range("N11").Value = 1.620655622136059
$ws.Range("O11").Value = 1.747323835194455

$ws.Range("I14").Value = 16.97685185185183
$ws.Range("N14").Value = 1.859533546038736
$ws.Range("O14").Value = 2.032558602498382

$ws.Range("I15").Value = -3.847222222222223
$ws.Range("N15").Value = 1.546865537736907
$ws.Range("O15").Value = 1.660778333536659

$ws.Range("I18").Value = 14.47727272727272
$ws.Range("N18").Value = 1.815485677363773
$ws.Range("O18").Value = 1.979371877230549

$ws.Range("I19").Value = 5.462962962962945
$ws.Range("N19").Value = 1.672603071948262
$ws.Range("O19").Value = 1.808689105403011

$ws.Range("I21").Value = 0.2777777777777778
$ws.Range("N21").Value = 1.600162412993039
$ws.Range("O21").Value = 1.723215189873418

$ws.Range("I22").Value = 1.791666666666668
$ws.Range("N22").Value = 1.620655622136059
$ws.Range("O22").Value = 1.747323835194455

$ws.Range("I23").Value = 13.46442495126706
$ws.Range("N23").Value = 1.798225615362447
$ws.Range("O23").Value = 1.958604378795604

$ws.Range("I24").Value = 19.79629629629628
$ws.Range("N24").Value = 1.911855479578636
$ws.Range("O24").Value = 2.09608909874769

$ws.Range("I28").Value = 21.28240740740739
$ws.Range("N28").Value = 1.940636870984383
$ws.Range("O28").Value = 2.131200751448103

$ws.Range("I29").Value = 1.791666666666668
$ws.Range("N29").Value = 1.620655622136059
$ws.Range("O29").Value = 1.747323835194455

$ws.Range("I31").Value = 12.67039049919483

$ws.Range("I32").Value = 12.67039049919483
$ws.Range("N32").Value = 1.784922174701128
$ws.Range("O32").Value = 1.942625691911729

$ws.Range("I33").Value = 19.65277777777778
$ws.Range("N33").Value = 1.909121107266436
$ws.Range("O33").Value = 2.092759415833974

$ws.Range("I34").Value = 21.28240740740739
$ws.Range("N34").Value = 1.940636870984383
$ws.Range("O34").Value = 2.131200751448103

$ws.Range("I35").Value = 13.0158303464755

$ws.Range("I36").Value = 13.0158303464755
$ws.Range("N36").Value = 1.790685487585954
$ws.Range("O36").Value = 1.94954496878686

$ws.Range("I37").Value = 19.60879629629628
$ws.Range("N37").Value = 1.908284719500103
$ws.Range("O37").Value = 2.091741145739967

$ws.Range("I39").Value = 14.47727272727272
$ws.Range("N39").Value = 1.815485677363773
$ws.Range("O39").Value = 1.979371877230549

$ws.Range("I40").Value = 14.47727272727272
$ws.Range("N40").Value = 1.815485677363773
$ws.Range("O40").Value = 1.979371877230549

$ws.Range("I41").Value = 14.96875
$ws.Range("N41").Value = 1.8239809580482
$ws.Range("O41").Value = 1.989608681354817

$ws.Range("I42").Value = 13.46442495126706

$ws.Range("I43").Value = 13.46442495126706
$ws.Range("N43").Value = 1.798225615362447
$ws.Range("O43").Value = 1.958604378795604

$ws.Range("I45").Value = 14.47727272727272
$ws.Range("N45").Value = 1.815485677363773
$ws.Range("O45").Value = 1.979371877230549

$ws.Range("I47").Value = 13.76976495726495
$ws.Range("N47").Value = 1.803394296576035
$ws.Range("O47").Value = 1.964819060413116

$ws.Range("I51").Value = 13.0158303464755
$ws.Range("N51").Value = 1.790685487585954
$ws.Range("O51").Value = 1.94954496878686

$ws.Range("I52").Value = 13.76976495726495
$ws.Range("N52").Value = 1.803394296576035
$ws.Range("O52").Value = 1.964819060413116

$ws.Range("I57").Value = 5.462962962962945
$ws.Range("N57").Value = 1.672603071948262
$ws.Range("O57").Value = 1.808689105403011
